# Refresh the crypto price/volume snapshot (Price=D, Volume(1h)=E).
# Values are stored as literal text (matching the sheet's existing
# inline-string cells), so a leading apostrophe forces text entry and
# ".Style = 'Normal'" strips the resulting quote-prefix formatting back
# to the original unstyled cell (no NumberFormat drift).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $text) {
    $rng = $ws.Range($cellRef)
    $rng.Value = "'" + $text
    $rng.Style = "Normal"
}

Set-TextValue "D2" "291.92"
Set-TextValue "E2" "0.24%"

Set-TextValue "D3" "30.97"
Set-TextValue "E3" "0.63%"

Set-TextValue "D4" "4.960"
Set-TextValue "E4" "1.26%"

Set-TextValue "D5" "0.07449"
Set-TextValue "E5" "2.68%"

Set-TextValue "D6" "2.243"
Set-TextValue "E6" "-4.73%"

Set-TextValue "D7" "7.740"
Set-TextValue "E7" "0.92%"

Set-TextValue "D8" "0.9191"
Set-TextValue "E8" "2.38%"

Set-TextValue "D9" "0.09358"
Set-TextValue "E9" "17.38%"

Set-TextValue "D10" "0.1724"
Set-TextValue "E10" "3.38%"

Set-TextValue "D11" "0.08363"
Set-TextValue "E11" "3.20%"

Set-TextValue "D12" "0.03252"
Set-TextValue "E12" "4.90%"

Set-TextValue "D13" "0.09932"
Set-TextValue "E13" "-0.80%"

Set-TextValue "D14" "0.001498"
Set-TextValue "E14" "-0.75%"

Set-TextValue "D15" "0.005768"
Set-TextValue "E15" "-1.19%"

Set-TextValue "D16" "3.470"

Set-TextValue "D17" "3.767"
Set-TextValue "E17" "1.66%"

Set-TextValue "E18" "3.47%"

Set-TextValue "D19" "0.3332"
Set-TextValue "E19" "0.40%"

Set-TextValue "D20" "0.1298"
Set-TextValue "E20" "0.02%"

Set-TextValue "D21" "4.141"
Set-TextValue "E21" "4.16%"

Set-TextValue "D22" "0.2123"
Set-TextValue "E22" "-7.88%"

Set-TextValue "D23" "0.04505"
Set-TextValue "E23" "-0.33%"

Set-TextValue "D24" "0.001218"
Set-TextValue "E24" "0.57%"

Set-TextValue "D25" "0.004259"
Set-TextValue "E25" "-3.21%"

Set-TextValue "D26" "0.0001297"
Set-TextValue "E26" "-0.31%"

Set-TextValue "D27" "0.0003393"
Set-TextValue "E27" "-0.16%"

Set-TextValue "D39" "0.01612"
Set-TextValue "E39" "2.18%"

Set-TextValue "D40" "0.04575"
Set-TextValue "E40" "4.57%"

Set-TextValue "D41" "0.007451"
Set-TextValue "E41" "1.83%"

Set-TextValue "D42" "0.009834"
Set-TextValue "E42" "-1.93%"

Set-TextValue "D43" "0.1357"
Set-TextValue "E43" "3.51%"

Set-TextValue "D44" "0.002155"
Set-TextValue "E44" "6.04%"

Set-TextValue "D45" "0.009191"
Set-TextValue "E45" "-3.29%"

Set-TextValue "D46" "0.00006100"
Set-TextValue "E46" "6.46%"

Set-TextValue "E47" "-0.28%"

Set-TextValue "D48" "2.525"
Set-TextValue "E48" "12.69%"

Set-TextValue "D49" "0.001999"
Set-TextValue "E49" "-31.12%"

Set-TextValue "E50" "-0.28%"

Set-TextValue "E51" "-0.28%"
